$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BP4D")
$ws.Range("B4").Value = 1
